$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9067
$ws.Range("C2").Value = 0.9067
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5798

$ws.Range("B3").Value = 0.992
$ws.Range("C3").Value = 0.992
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.8033

$ws.Range("B4").Value = 0.9896
$ws.Range("C4").Value = 0.9896
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.8442

$ws.Range("B5").Value = 0.9874000000000001
$ws.Range("C5").Value = 0.9874000000000001
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.8874

$ws.Range("B6").Value = 0.9791
$ws.Range("C6").Value = 0.9791
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.8509

$ws.Range("B7").Value = 0.9715
$ws.Range("C7").Value = 0.9715
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.8243

$ws.Range("B8").Value = 0.9471000000000001
$ws.Range("C8").Value = 0.9471000000000001
$ws.Range("D8").Value = 0.9494
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.8022
